# Update crypto price (column D) and 1h volume/change (column E) values
# to match the refreshed data from the coinranking.com feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.640.33"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.578.34"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.52"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.82"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.04"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0592"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "1.573.20"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "28.664.50"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.35"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.24"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.89"
$ws.Range("E23").Value = "  -4.18%  "
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +6.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.71"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.03"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.46"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "1.399.00"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.524"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0465"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.961"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.09"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "1.715.78"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.59"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -0.94%  "
